$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$headers = @(
    "Inhaber",
    "total_km_durchgangsstrasse",
    "kb_befreit",
    "AmpelcodePers1",
    "AmpelcodePers2",
    "AmpelcodePers3",
    "AmpelcodePers4",
    "AmpelcodePers5",
    "AmpelcodeOFG1",
    "AmpelcodeOFG2",
    "AmpelcodeOFG3",
    "AmpelcodeOFG5",
    "AmpelcodeGW1",
    "AmpelcodeGW3",
    "AmpelcodeGW4",
    "AmpelcodeGW5"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Build the header style (bold font, thin box border, centered/top aligned)
# on a single cell first, so the engine only mints one new font / border /
# cellXf each - then fan it out to the rest of the header row by copying
# the computed format, instead of re-deriving the style per cell.
$firstHeaderCell = $ws.Range("A1")
$firstHeaderCell.HorizontalAlignment = -4108   # xlCenter
$firstHeaderCell.VerticalAlignment = -4160     # xlTop
$firstHeaderCell.Borders.LineStyle = 1         # xlContinuous
$firstHeaderCell.Font.Bold = $true

$firstHeaderCell.Copy()
$ws.Range("B1:P1").PasteSpecial(-4122)         # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows ---
$data = @(
    @("LU", 365.898, 0, 346.457, 8.952, 7.513, 2.971, 0.005, 343.37, 8.132999999999999, 14.39, 0.005, 344.948, 4.605, 16.34, 0.005),
    @("TG", 344.233, 0, 333.885, 6.685, 1.66, 0.54, 1.463, 335.28, 5.345, 2.145, 1.463, 284.425, 2.928, 55.417, 1.463),
    @("BL", 176.2151, 0, 162.7862, 13.4289, $null, $null, $null, 131.8565, 26.0675, 18.2911, $null, 169.2787, 2.9284, 4.008, $null)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        if ($null -ne $value) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $value
        }
    }
}
